# Cat Wilde and the Doom of Dead - add meta description paragraph near the
# top, drop the duplicated "Play Cat Wilde..." heading-style paragraph near
# the bottom, and rewrite the final (italic) paragraph with the image prompt.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph: "Play Cat Wilde and the Doom of Dead Free |
#    Review". New paragraph = Normal style, with a bold "Meta description"
#    lead-in followed by plain text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start
$metaFull = "Meta description: Experience ancient Egypt with Cat Wilde and the Doom of Dead slot game by Play N Go for free. Read our review for game features, design, and gameplay."
$metaRange.Text = $metaFull

$boldLabel = "Meta description"
$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Near the end of the document, remove the paragraph that duplicates
#    the bold "Play Cat Wilde and the Doom of Dead Free | Review" text,
#    and replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping its italic formatting intact.
#    (Search starts after the real Heading1 title so the first,
#    legitimate occurrence of the title text is left untouched.)
# ---------------------------------------------------------------------
$searchStart = $titleRange.End
$dupRange = $d.Range($searchStart, $d.Content.End)
$found = $dupRange.Find.Execute("Play Cat Wilde and the Doom of Dead Free | Review")

if ($found) {
    [void]$dupRange.Expand(4)
    $dupRange.Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastStart = $lastRange.Start
$lastEnd = $lastRange.End - 1
$lastTextRange = $d.Range($lastStart, $lastEnd)

$promptText = "Prompt: Create a cartoon-style feature image for " + [char]34 + "Cat Wilde and the Doom of Dead" + [char]34 + " slot game with a happy Maya warrior wearing glasses as the main focus. The image should be visually stunning with bold colors and a playful, adventurous tone. In the center of the image, feature Cat Wilde, the adventurous explorer of ancient Egypt, standing confident and smiling with a fierce warrior outfit. Add some sand dunes and pyramids in the background, and have the Eye of Ra symbol, which acts as the wild in the game, hovering above Wilde's head and shining brightly. To Wilde's right, draw a happy Maya warrior wearing glasses and holding a compass and a sarcophagus. Make sure the warrior is looking pleased and excited about the treasure hunting adventure. In the top left corner of the image, add the game's title, " + [char]34 + "Cat Wilde and the Doom of Dead" + [char]34 + " in bold letters with the subtitle " + [char]34 + "Join Cat Wilde on an ancient Egyptian adventure" + [char]34 + " just below. Make sure the image is visually appealing, age-appropriate, and attention-grabbing enough to entice players to give the game a try."

$lastTextRange.Text = $promptText
